$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 7 ---
$ws.Range("A7").Value = "Date"
$ws.Range("B7").Value = "Location"
$ws.Range("C7").Value = "Works"
$ws.Range("D7").Value = "Workers"
$ws.Range("E7").Value = "Hours"
$ws.Range("F7").Value = "Project amounts"
$ws.Range("G7").Value = "Completed "

# --- Row 8 ---
$ws.Range("B8").Value = "Fondations"
$ws.Range("C8").Value = "Concreeting "
$ws.Range("D8").Value = 3
$ws.Range("E8").Formula = "=D8*4"
$ws.Range("F8").Value = "80m3"
$ws.Range("G8").Value = "20m3"

# --- Row 9 ---
$ws.Range("B9").Value = "Fondations"
$ws.Range("C9").Value = "Finishin works "
$ws.Range("D9").Value = 2
$ws.Range("E9").Value = 15
$ws.Range("F9").Value = "500m2"
$ws.Range("G9").Value = "100m2"

# --- Row 10 ---
$ws.Range("B10").Value = "First floor"
$ws.Range("C10").Value = "Assembly"
$ws.Range("D10").Value = 4
$ws.Range("E10").Formula = "=4*9"
$ws.Range("F10").Value = 60
$ws.Range("G10").Value = 25

# --- Section headers ---
$ws.Range("A12").Value = "Delivered materials"
$ws.Range("A17").Value = "Photos : "

# --- Outer box borders (medium) around each of the three tables ---
$ws.Range("A7:G10").BorderAround(1, -4138)
$ws.Range("A12:G15").BorderAround(1, -4138)
$ws.Range("A17:G21").BorderAround(1, -4138)

# --- Materialise the untouched interior cells of the blank template rows
#     (they carry a "no border applied" style in the authored sheet) ---
$ws.Range("B13:F14").Borders.Item(7).LineStyle = -4142
$ws.Range("B18:F20").Borders.Item(7).LineStyle = -4142

# --- Right-align the "Completed" amount cells ---
$ws.Range("F8:F9").HorizontalAlignment = -4152
$ws.Range("G8:G10").HorizontalAlignment = -4152

# --- Row heights for the thick-bottom-border separator / closing rows ---
$ws.Rows.Item(6).RowHeight = 15.75
$ws.Rows.Item(10).RowHeight = 15.75
$ws.Rows.Item(11).RowHeight = 15.75
$ws.Rows.Item(15).RowHeight = 15.75
$ws.Rows.Item(16).RowHeight = 15.75
$ws.Rows.Item(21).RowHeight = 15.75

# --- Column widths (closest achievable to the authored pixel widths) ---
$ws.Columns.Item(1).ColumnWidth = 17.8333333333333
$ws.Columns.Item(2).ColumnWidth = 13.6666666666667
$ws.Columns.Item(3).ColumnWidth = 13.6666666666667
$ws.Columns.Item(6).ColumnWidth = 16.5
$ws.Columns.Item(7).ColumnWidth = 11.5

# --- Selection / active cell matching the authored view ---
$ws.Range("A7:G21").Select()
